$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.463.35"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.891.21"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.60"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "0.4844"
$ws.Range("E7").Value = "  -1.50%  "
$ws.Range("D8").Value = "0.2892"
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("D9").Value = "0.06631"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").Value = "1.906.65"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").Value = "16.98"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "0.07413"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "5.202"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "89.15"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6650"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "30.429.87"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "13.55"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007788"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").Value = "0.9988"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "5.434"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("D21").Value = "2.128.80"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").Value = "0.9994"
$ws.Range("D23").Value = "217.55"
$ws.Range("E23").Value = "  +11.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.210"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.450"
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("D26").Value = "165.18"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("D27").Value = "18.65"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").Value = "1.947"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.440"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("D30").Value = "4.323"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09200"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.080"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "0.05085"
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("D34").Value = "0.7512"
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("D35").Value = "1.158"
$ws.Range("E35").Value = "  +4.52%  "
$ws.Range("D36").Value = "2.701"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").Value = "0.01891"
$ws.Range("E37").Value = "  +4.14%  "
$ws.Range("D38").Value = "2.649"
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9210"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").Value = "2.095"
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("D41").Value = "6.099"
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("D42").Value = "107.52"
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "7.654"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").Value = "0.1351"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "66.01"
$ws.Range("E47").Value = "  -11.76%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.575"
$ws.Range("E48").Value = "  +10.40%  "
$ws.Range("D49").Value = "8.944"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").Value = "34.43"
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").Value = "0.05699"
$ws.Range("E51").Value = "  -2.75%  "
